# FarmerTable: insert a new "rarity" column between the existing id/nameLocalKey
# columns, carrying the same per-row formatting as the column it pushes aside,
# then populate the new column's header/type/value triplet (rarity / Erarity / Common).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "nameLocalKey" column (B) one slot to the right (-> C),
# opening up a blank column B for the new field.
$ws.Range("B2:B4").Insert(-4161)

# The freshly-opened column B is currently styled like column A (Insert()
# copies the format of the column being pushed away from). Re-stamp each row
# of column B with the formatting of column C instead, since the new
# "rarity" field is conceptually a header/type/value column just like the
# nameLocalKey column it now sits beside.
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new field definition: name / type / default value.
$ws.Range("B2").Value = "rarity"
$ws.Range("B3").Value = "Erarity"
$ws.Range("B4").Value = "Common"

# Approximate column widths close to the post-edit layout.
$ws.Columns("B").ColumnWidth = 12.86

# Restore the cursor to where the editing user last left it.
[void]$ws.Range("F10").Select()
